# إضافة عمود جديد 'Event ' إلى Card21
# Also normalizes Card22: literal "nan" placeholder text is cleared out and the
# trailing blank row (row 13) is removed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Card22: clear every cell that currently holds the literal text "nan" so the
# cell becomes blank, then delete the now-unused trailing blank row 13.
# ---------------------------------------------------------------------------
$ws22 = $wb.Worksheets.Item("Card22")

$lastRow22 = 13
$lastCol22 = 15  # column O

for ($r = 2; $r -le $lastRow22; $r++) {
    for ($c = 1; $c -le $lastCol22; $c++) {
        $cell = $ws22.Cells.Item($r, $c)
        $v = $cell.Value()
        if (-not [string]::IsNullOrEmpty($v)) {
            if ($v.Equals("nan")) {
                $cell.Value = ""
            }
        }
    }
}

$ws22.Rows.Item(13).Delete()

# ---------------------------------------------------------------------------
# Card21: fill every currently blank data cell (columns D..L, rows 2..12)
# with the literal text "nan" to match the other cards, then append the new
# 'Event ' column (M) with the same header styling as the other headers.
# ---------------------------------------------------------------------------
$ws21 = $wb.Worksheets.Item("Card21")

$firstDataCol21 = 4   # column D
$lastDataCol21 = 12   # column L
$lastRow21 = 12

for ($r = 2; $r -le $lastRow21; $r++) {
    for ($c = $firstDataCol21; $c -le $lastDataCol21; $c++) {
        $cell = $ws21.Cells.Item($r, $c)
        $v = $cell.Value()
        if ([string]::IsNullOrEmpty($v)) {
            $cell.Value = "nan"
        }
    }
}

# New header column M: "Event " (note trailing space), styled like the other
# header cells (bold font, thin border, centered/top aligned).
$ws21.Range("L1").Copy()
$ws21.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws21.Range("M1").Value = "Event "
